# Applies the cryptos-list refresh described by the commit:
# "Updated cryptos list on Fri Jul 19 13:32:34 UTC 2024 with GitHub Actions"
#
# Strategy: write each changed cell's new text back with .Value. Column D
# holds text in the source workbook (t="inlineStr"), but some of the new
# values look like plain numbers (e.g. "574.44") and Excel's COM layer would
# silently coerce those to numeric cells. To keep them as text - matching the
# original inline-string cells - we set NumberFormat = "@" (Text) on exactly
# those cells right before assigning the value. Values that are unambiguous
# text (URLs, names, percentages with spaces, multi-dot price strings like
# "64.144.57") are left on the default format since Excel keeps them as text
# anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.144.57"
$ws.Range("E2").Value = "  -1.13%  "
# Row 3
$ws.Range("D3").Value = "3.403.18"
$ws.Range("E3").Value = "  -2.22%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.44"
$ws.Range("E5").Value = "  -0.48%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.71"
$ws.Range("E6").Value = "  +0.85%  "
# Row 7
$ws.Range("E7").Value = "  +0.02%  "
# Row 8
$ws.Range("D8").Value = "3.402.72"
$ws.Range("E8").Value = "  -2.25%  "
# Row 9
$ws.Range("E9").Value = "  -4.89%  "
# Row 10
$ws.Range("E10").Value = "  +0.74%  "
# Row 11
$ws.Range("E11").Value = "  -3.06%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.419"
$ws.Range("E12").Value = "  -4.54%  "
# Row 13
$ws.Range("D13").Value = "3.985.76"
$ws.Range("E13").Value = "  -2.30%  "
# Row 14
$ws.Range("E14").Value = "  +0.16%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.79"
$ws.Range("E15").Value = "  -3.15%  "
# Row 16
$ws.Range("E16").Value = "  -2.95%  "
# Row 17
$ws.Range("D17").Value = "64.141.80"
$ws.Range("E17").Value = "  -1.24%  "
# Row 18
$ws.Range("D18").Value = "3.392.76"
$ws.Range("E18").Value = "  -4.46%  "
# Row 19
$ws.Range("E19").Value = "  -1.96%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.43"
$ws.Range("E20").Value = "  -3.20%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.10"
$ws.Range("E21").Value = "  -2.42%  "
# Row 22
$ws.Range("E22").Value = "  -2.89%  "
# Row 23
$ws.Range("E23").Value = "  -0.01%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.33"
# Row 25
$ws.Range("E25").Value = "  -3.88%  "
# Row 26
$ws.Range("E26").Value = "  -5.59%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.43"
$ws.Range("E27").Value = "  -4.94%  "
# Row 28
$ws.Range("E28").Value = "  -0.73%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.13%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.05"
$ws.Range("E30").Value = "  -1.46%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.39"
$ws.Range("E31").Value = "  -4.23%  "
# Row 32
$ws.Range("E32").Value = "  -0.96%  "
# Row 33
$ws.Range("E33").Value = "  +0.02%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.77"
$ws.Range("E34").Value = "  -2.85%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.99"
$ws.Range("E35").Value = "  -1.39%  "
# Row 36
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.55"
$ws.Range("E36").Value = "  -1.22%  "
# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.47"
$ws.Range("E37").Value = "  -8.03%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.860"
$ws.Range("E38").Value = "  +5.37%  "
# Row 39
$ws.Range("E39").Value = "  -3.84%  "
# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0719"
$ws.Range("E40").Value = "  -4.96%  "
# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.759.70"
$ws.Range("E41").Value = "  -4.25%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "25.58"
$ws.Range("E42").Value = "  -5.18%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.50"
$ws.Range("E43").Value = "  -1.20%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.41"
$ws.Range("E44").Value = "  -2.44%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.74"
$ws.Range("E45").Value = "  -1.33%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.33"
$ws.Range("E46").Value = "  -4.59%  "
# Row 47
$ws.Range("E47").Value = "  -2.69%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  -1.68%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "327.48"
$ws.Range("E49").Value = "  -1.43%  "
# Row 50
$ws.Range("E50").Value = "  -3.26%  "
# Row 51
$ws.Range("E51").Value = "  -3.54%  "
